$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text (string) cell type for the Price/Volume columns while writing
# numeric-looking literals, then restore the original (unstyled) cell style so
# no incidental formatting diff is introduced.
$fmtRange = $ws.Range("D2:E51")
$fmtRange.NumberFormat = "@"

$ws.Range("D2").Value = "304.35"
$ws.Range("E2").Value = "3.91%"

$ws.Range("D3").Value = "35.72"
$ws.Range("E3").Value = "14.44%"

$ws.Range("D4").Value = "5.089"
$ws.Range("E4").Value = "2.46%"

$ws.Range("D5").Value = "0.07834"
$ws.Range("E5").Value = "4.55%"

$ws.Range("D6").Value = "2.246"
$ws.Range("E6").Value = "-0.86%"

$ws.Range("D7").Value = "8.125"
$ws.Range("E7").Value = "4.28%"

$ws.Range("D8").Value = "4.002"
$ws.Range("E8").Value = "5.87%"

$ws.Range("D9").Value = "0.9260"
$ws.Range("E9").Value = "0.75%"

$ws.Range("D10").Value = "0.09861"
$ws.Range("E10").Value = "6.82%"

$ws.Range("D11").Value = "0.1818"
$ws.Range("E11").Value = "5.09%"

$ws.Range("D12").Value = "0.08682"
$ws.Range("E12").Value = "4.12%"

$ws.Range("D13").Value = "0.03413"
$ws.Range("E13").Value = "4.04%"

$ws.Range("D14").Value = "0.09926"
$ws.Range("E14").Value = "-0.13%"

$ws.Range("D15").Value = "0.001471"
$ws.Range("E15").Value = "-1.83%"

$ws.Range("D16").Value = "0.005775"
$ws.Range("E16").Value = "1.07%"

$ws.Range("D17").Value = "3.488"
$ws.Range("E17").Value = "0.38%"

$ws.Range("E18").Value = "-1.24%"

$ws.Range("E19").Value = "2.51%"

$ws.Range("E20").Value = "0.40%"

$ws.Range("D21").Value = "4.548"
$ws.Range("E21").Value = "11.33%"

$ws.Range("D23").Value = "0.04678"
$ws.Range("E23").Value = "3.11%"

$ws.Range("D24").Value = "0.001239"
$ws.Range("E24").Value = "1.60%"

$ws.Range("D25").Value = "0.004511"
$ws.Range("E25").Value = "4.78%"

$ws.Range("D26").Value = "0.0001302"
$ws.Range("E26").Value = "0.34%"

$ws.Range("D27").Value = "0.0002697"
$ws.Range("E27").Value = "-20.55%"

$ws.Range("D39").Value = "0.01762"
$ws.Range("E39").Value = "8.41%"

$ws.Range("D40").Value = "0.04706"
$ws.Range("E40").Value = "2.69%"

$ws.Range("D41").Value = "0.008028"
$ws.Range("E41").Value = "7.26%"

$ws.Range("D42").Value = "0.1421"
$ws.Range("E42").Value = "4.47%"

$ws.Range("D43").Value = "0.008470"
$ws.Range("E43").Value = "-13.95%"

$ws.Range("E44").Value = "-0.11%"

$ws.Range("D45").Value = "0.009127"
$ws.Range("E45").Value = "-6.71%"

$ws.Range("D46").Value = "0.00006160"
$ws.Range("E46").Value = "0.82%"

$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").Value = "-0.08%"

$ws.Range("D48").Value = "5.742"
$ws.Range("E48").Value = "105.49%"

$ws.Range("D49").Value = "0.002688"
$ws.Range("E49").Value = "34.38%"

$ws.Range("D50").Value = "0.00002098"
$ws.Range("E50").Value = "-0.08%"

$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").Value = "-0.08%"

$fmtRange.Style = "Normal"
